$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.316.71"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").Value = "1.588.83"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  -0.46%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.07"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.41"
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0845"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").Value = "1.812.17"
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").Value = "1.564.71"
$ws.Range("E14").Value = "  -1.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.520"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.38"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "26.322.16"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.45"
$ws.Range("E19").Value = "  +5.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "210.91"
$ws.Range("E20").Value = "  +1.53%  "
$ws.Range("E21").Value = "  -0.43%  "
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("E24").Value = "  -2.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.61"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.24"
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("E32").Value = "  -0.87%  "
$ws.Range("E33").Value = "  +1.43%  "
$ws.Range("D34").Value = "1.316.11"
$ws.Range("E34").Value = "  +2.60%  "
$ws.Range("E36").Value = "  +1.73%  "
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("E39").Value = "  -13.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.805"
$ws.Range("E40").Value = "  -1.66%  "
$ws.Range("E41").Value = "  -0.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.64"
$ws.Range("E42").Value = "  +3.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.767"
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("E44").Value = "  -1.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.28"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("D46").Value = "1.725.12"
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.51"
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("E48").Value = "  -5.42%  "
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("E50").Value = "  -4.83%  "
$ws.Range("E51").Value = "  -0.41%  "
